$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column width changes (raw OOXML width = ColumnWidth + 0.83) ---
$ws.Columns.Item(2).ColumnWidth = 84.17
$ws.Columns.Item(3).ColumnWidth = 84.17

$ws.Columns.Item(21).ColumnWidth = 34.17
$ws.Columns.Item(22).ColumnWidth = 36.17
$ws.Columns.Item(23).ColumnWidth = 41.17
$ws.Columns.Item(24).ColumnWidth = 43.17
$ws.Columns.Item(25).ColumnWidth = 31.17
$ws.Columns.Item(26).ColumnWidth = 33.17
$ws.Columns.Item(27).ColumnWidth = 28.17
$ws.Columns.Item(28).ColumnWidth = 30.17
$ws.Columns.Item(29).ColumnWidth = 37.17
$ws.Columns.Item(30).ColumnWidth = 39.17
$ws.Columns.Item(31).ColumnWidth = 23.17
$ws.Columns.Item(32).ColumnWidth = 25.17
$ws.Columns.Item(33).ColumnWidth = 24.17
$ws.Columns.Item(34).ColumnWidth = 26.17
$ws.Columns.Item(35).ColumnWidth = 27.17
$ws.Columns.Item(36).ColumnWidth = 29.17
$ws.Columns.Item(37).ColumnWidth = 24.17
$ws.Columns.Item(38).ColumnWidth = 26.17

$ws.Columns.Item(50).ColumnWidth = 25.17

# --- Row 1 header renames ---
$ws.Range("A1").Value = "button_closeActions_class"

$ws.Range("U1").Value = "link_executionLinks_executions_id"
$ws.Range("V1").Value = "link_executionLinks_executions_id_1"
$ws.Range("W1").Value = "link_executionLinks_internalRoleLinkName"
$ws.Range("X1").Value = "link_executionLinks_internalRoleLinkName_1"
$ws.Range("Y1").Value = "link_executionLinks_project_id"
$ws.Range("Z1").Value = "link_executionLinks_project_id_1"
$ws.Range("AA1").Value = "link_executionLinks_team_id"
$ws.Range("AB1").Value = "link_executionLinks_team_id_1"
$ws.Range("AC1").Value = "link_jobDetails_internalRoleLinkName"
$ws.Range("AD1").Value = "link_jobDetails_internalRoleLinkName_1"
$ws.Range("AE1").Value = "link_jobDetails_job_id"
$ws.Range("AF1").Value = "link_jobDetails_job_id_1"
$ws.Range("AG1").Value = "link_jobDetails_plan_id"
$ws.Range("AH1").Value = "link_jobDetails_plan_id_1"
$ws.Range("AI1").Value = "link_jobDetails_project_id"
$ws.Range("AJ1").Value = "link_jobDetails_project_id_1"
$ws.Range("AK1").Value = "link_jobDetails_team_id"
$ws.Range("AL1").Value = "link_jobDetails_team_id_1"

$ws.Range("AX1").Value = "span_logContent_nthChild"

# --- Row 2 data changes ---
$ws.Range("B2").Value = "Data Files/AI-Generated/Common/scheduleAndRunTestWithEnvironmentSelection-test-data"
$ws.Range("C2").Value = "Data Files/AI-Generated/Common/scheduleAndRunTestWithEnvironmentSelection-test-data"

# Leading apostrophe forces these (numeric-looking) entries to stay text,
# matching the original data's text storage instead of being coerced to numbers.
$ws.Range("U2").Value = "'10"
$ws.Range("V2").Value = "'12"
$ws.Range("W2").Value = "'10"
$ws.Range("X2").Value = "'12"
$ws.Range("Y2").Value = "'1588984"
$ws.Range("Z2").Value = "'1588984"
$ws.Range("AA2").Value = "'1570311"
$ws.Range("AB2").Value = "'1570311"
$ws.Range("AC2").Value = "'8"
$ws.Range("AD2").Value = "'10"
$ws.Range("AE2").Value = "'8"
$ws.Range("AF2").Value = "'10"
$ws.Range("AG2").Value = "'837097"
$ws.Range("AH2").Value = "'837132"
